$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Édité le / Par :" header line with the new timestamp and the
#    user name who is logging on (the actual edit requested by the commit).
$ws.Range("A1").Value = "Édité le : 31/01/2025 à 23:56:46`nPar : FouedAmich"
# Re-fit the row height back to the sheet default; Excel auto-expands the row
# for the embedded newline when the cell value is assigned, which would add
# an explicit ht/customHeight that isn't part of the intended change.
$ws.Rows.Item(1).AutoFit()

# 2. Clear out the old "Le : dd/mm/yyyy" line.
$ws.Range("A3").Value = ""

# 3. Drop the now-unused trailing column F from the report grid.
$ws.Columns.Item(6).Delete()

# Deleting the column collapses the trivial single-cell merges on A6/A7;
# restore them so the merge map matches the original layout.
$ws.Range("A6:A6").Merge()
$ws.Range("A7:A7").Merge()

# 4. Even out the remaining columns (A:E) to a uniform width.
$w = 16 - (5/6)
$ws.Columns.Item(1).ColumnWidth = $w
$ws.Columns.Item(2).ColumnWidth = $w
$ws.Columns.Item(3).ColumnWidth = $w
$ws.Columns.Item(4).ColumnWidth = $w
$ws.Columns.Item(5).ColumnWidth = $w
